# Applies the cryptos price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.809.47"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.800.35"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.55"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.47"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.88"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "4.436.33"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "3.816.68"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "67.834.32"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.38"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.39"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.88"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.05"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.11"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D29").Value = "3.947.94"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.46"
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.78"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.988"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.77"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.58"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "150.78"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("E47").Value = "  +8.78%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "399.58"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.40"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.86"
$ws.Range("E51").Value = "  +1.93%  "
